$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
try {
  $st = $wb.Styles.Add("HeaderLR")
  Write-Host "Added style"
  $st.Font.Bold = $true
  $st.Borders.Item(7).LineStyle = 1
  $st.Borders.Item(10).LineStyle = 1
  $ws.Range("K1").Style = "HeaderLR"
  Write-Host "done"
} catch {
  Write-Host "ERR:" $_
}
